{"js": "const body = context.document.body;\nconst replacements = [\n  { old: \"70\u00d757=3990\", news: [\"29\u00d770=2030\"] },\n  { old: \"13\u00d786=1118\", news: [\"80\u00d798=7840\"] },\n  { old: \"68\u00d719=1292\", news: [\"25\u00d720=500\"] },\n  { old: \"45\u00d759=2655\", news: [\"17\u00d787=1479\"] },\n  { old: \"44\u00d767=2948\", news: [\"32\u00d788=2816\"] },\n  { old: \"82\u00d769=5658\", news: [\"12\u00d726=312\"] },\n  { old: \"87\u00d729=2523\", news: [\"71\u00d784=5964\"] },\n  { old: \"36\u00d776=2736\", news: [\"44\u00d730=1320\", \"20\u00d727=540\"] },\n  { old: \"54\u00d716=864\", news: [\"84\u00d731=2604\"] },\n  { old: \"98\u00d797=9506\", news: [\"69\u00d779=5451\"] },\n  { old: \"99\u00d740=3960\", news: [\"45\u00d796=4320\"] },\n  { old: \"35\u00d750=1750\", news: [\"18\u00d721=378\"] },\n  { old: \"90\u00d746=4140\", news: [\"96\u00d757=5472\"] },\n  { old: \"82\u00d762=5084\", news: [\"85\u00d737=3145\"] },\n  { old: \"30\u00d764=1920\", news: [\"62\u00d752=3224\"] },\n  { old: \"47\u00d794=4418\", news: [\"58\u00d772=4176\"] },\n  { old: \"11\u00d799=1089\", news: [\"98\u00d745=4410\"] },\n  { old: \"90\u00d741=3690\", news: [\"45\u00d785=3825\"] },\n  { old: \"43\u00d750=2150\", news: [\"73\u00d746=3358\"] },\n  { old: \"27\u00d766=1782\", news: [\"60\u00d765=3900\"] },\n  { old: \"50\u00d716=800\", news: [\"97\u00d725=2425\"] },\n  { old: \"32\u00d785=2720\", news: [\"84\u00d736=3024\"] },\n  { old: \"36\u00d738=1368\", news: [\"21\u00d750=1050\"] },\n  { old: \"35\u00d760=2100\", news: [\"46\u00d772=3312\"] },\n];\n\nfor (const { old, news } of replacements) {\n  const results = body.search(old, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== news.length) {\n    throw new Error(`Expected ${news.length} matches for \"${old}\", found ${results.items.length}`);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(news[i], Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"70\u00d757=3990\"; New = @(\"29\u00d770=2030\") },\n    @{ Old = \"13\u00d786=1118\"; New = @(\"80\u00d798=7840\") },\n    @{ Old = \"68\u00d719=1292\"; New = @(\"25\u00d720=500\") },\n    @{ Old = \"45\u00d759=2655\"; New = @(\"17\u00d787=1479\") },\n    @{ Old = \"44\u00d767=2948\"; New = @(\"32\u00d788=2816\") },\n    @{ Old = \"82\u00d769=5658\"; New = @(\"12\u00d726=312\") },\n    @{ Old = \"87\u00d729=2523\"; New = @(\"71\u00d784=5964\") },\n    @{ Old = \"36\u00d776=2736\"; New = @(\"44\u00d730=1320\", \"20\u00d727=540\") },\n    @{ Old = \"54\u00d716=864\"; New = @(\"84\u00d731=2604\") },\n    @{ Old = \"98\u00d797=9506\"; New = @(\"69\u00d779=5451\") },\n    @{ Old = \"99\u00d740=3960\"; New = @(\"45\u00d796=4320\") },\n    @{ Old = \"35\u00d750=1750\"; New = @(\"18\u00d721=378\") },\n    @{ Old = \"90\u00d746=4140\"; New = @(\"96\u00d757=5472\") },\n    @{ Old = \"82\u00d762=5084\"; New = @(\"85\u00d737=3145\") },\n    @{ Old = \"30\u00d764=1920\"; New = @(\"62\u00d752=3224\") },\n    @{ Old = \"47\u00d794=4418\"; New = @(\"58\u00d772=4176\") },\n    @{ Old = \"11\u00d799=1089\"; New = @(\"98\u00d745=4410\") },\n    @{ Old = \"90\u00d741=3690\"; New = @(\"45\u00d785=3825\") },\n    @{ Old = \"43\u00d750=2150\"; New = @(\"73\u00d746=3358\") },\n    @{ Old = \"27\u00d766=1782\"; New = @(\"60\u00d765=3900\") },\n    @{ Old = \"50\u00d716=800\"; New = @(\"97\u00d725=2425\") },\n    @{ Old = \"32\u00d785=2720\"; New = @(\"84\u00d736=3024\") },\n    @{ Old = \"36\u00d738=1368\"; New = @(\"21\u00d750=1050\") },\n    @{ Old = \"35\u00d760=2100\"; New = @(\"46\u00d772=3312\") },\n)\n\nforeach ($item in $replacements) {\n    $range = $d.Content\n    foreach ($newText in $item.New) {\n        $range.Find.ClearFormatting()\n        $range.Find.Text = $item.Old\n        $range.Find.Forward = $true\n        $range.Find.Wrap = 0\n        $found = $range.Find.Execute()\n        if (-not $found) {\n            throw \"Could not find occurrence of $($item.Old)\"\n        }\n        $range.Text = $newText\n        $range.Collapse(0)\n    }\n}"}
